$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("A2:I36")
$key1 = $ws.Range("B2:B36")
$key2 = $ws.Range("A2:A36")
$rng.Sort($key1, 1, $key2, 0, 1)
